# Update RPA test cache workbook: correct/garble a batch of menu-item
# description strings (shared strings) on both the "Items" and the
# "Items - Formatted" worksheets. Column A, rows 3-12 hold the item
# descriptions affected by this change.

$wb = $excel.ActiveWorkbook

$updates = @{
    "A3"  = "Pan Fried Leek Dumplings IAT (2)"
    "A4"  = "Pork Xiao Long Bao(10) A¥R]J¿E(10)"
    "A5"  = "Q-BAO (5) WEEL (5)"
    "A6"  = "Chicken potstickers KR`$55(6)"
    "A7"  = "Tomato Mushroom Steamed dumpli pEiAINABUNXA (6)"
    "A8"  = "Zucchini shrimp dumplings A/LC"
    "A9"  = "beef stew nodle soup (Non Spicy 0H#4PJB(TY)"
    "A10" = "dandan noodle INCMM"
    "A11" = "banana naan bread BATAI"
    "A12" = "house made plum juice"
}

foreach ($sheetName in @("Items", "Items - Formatted")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($addr in $updates.Keys) {
        $ws.Range($addr).Value = $updates[$addr]
    }
}
